$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the page names for the Subpopulation and LOT pages (row 6)
$ws.Range("F6").Value = "pop_filter1_section1"
$ws.Range("G6").Value = "pop_filter1_section1_checkbox"
$ws.Range("H6").Value = "pop_filter1_section"

# Move the active selection to F10, matching the author's final cursor position
$ws.Range("F10").Select()
